# ----------------------------------------------------------------------
# Rapport SEO.xlsx — rework the "Liste" task-template sheet into a two
# table SEO / Performance audit sheet (per commit "Modification des
# fichiers de rapport selon les directives du projet").
# ----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1. Title block (row 1 / row 2, still merged A1:C1 / A2:C2) -------
$ws.Range("A1").Value2 = "Analyse SEO / Performance"
$ws.Range("A2").Value2 = "Audit"

# ---- 2. New column widths for the added columns D, E, F ---------------
$ws.Columns.Item(4).ColumnWidth = 30.140625
$ws.Columns.Item(5).ColumnWidth = 27.42578125
$ws.Columns.Item(6).ColumnWidth = 27.42578125

# ---- 3. Row 3 becomes the (only) header row for both tables -----------
# Left table ("Liste", A3:C7) header text
$ws.Range("B3").Value2 = "Problème identifié"
$ws.Range("C3").Value2 = "Explication du problème"
# Right table ("Liste3", D3:F7) header text
$ws.Range("D3").Value2 = "Bonne pratique à adopter"
$ws.Range("E3").Value2 = "Action recommandée"
$ws.Range("F3").Value2 = "Référence"
$ws.Range("A3").Value2 = "Catégorie (SEO? Performance? Accessibilité?)"

# ---- 4. Data rows 4-7: mirror the Date / Élément / Remarque sample ----
#         formatting from A:C into the new D:F columns, then copy over --
#         the same sample text values --------------------------------
$ws.Range("A4:C7").Copy()
$ws.Range("D4").PasteSpecial(-4122)   # xlPasteFormats
for ($r = 4; $r -le 7; $r++) {
    $ws.Cells.Item($r, 4).Value2 = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 5).Value2 = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 6).Value2 = $ws.Cells.Item($r, 3).Value2
}

# ---- 5. Row heights: header row 3 grows to match rows 4-7 (30pt) ------
$ws.Rows.Item(3).RowHeight = 30

# ---- 6. Second table "Liste3" over D3:F7 (created before the manual ---
#         header paint job below so both tables' header rows pick up ----
#         the same explicit formatting, like the real workbook) --------
$table2 = $ws.ListObjects.Add(1, $ws.Range("D3:F7"), 0, 1)
$table2.Name = "Liste3"
$table2.TableStyle = "Tâches"
$table2.ShowTableStyleFirstColumn = $true

# ---- 7. Header-row look: solid blue fill + white wrapped text ---------
$headerRange = $ws.Range("A3:F3")
$headerRange.Interior.Color = 15773696          # RGB(0,176,240) = FF00B0F0
$headerRange.Font.Color = 16777215              # white
$headerRange.VerticalAlignment = -4108          # xlVAlignCenter
$headerRange.WrapText = $true

# ---- 8. Data-validation prompts on the new header cells (mirrors the --
#         ones already present on A3 / B3 / C3) -------------------------
$ws.Range("D3").Validation.Add(0, 1, 1, "")
$ws.Range("D3").Validation.InputMessage = "Entrez la date dans cette colonne sous ce titre. Utilisez les filtres des titres pour trouver des entrées spécifiques"
$ws.Range("D3").Validation.ShowInput = $true
$ws.Range("D3").Validation.ShowError = $true

$ws.Range("E3").Validation.Add(0, 1, 1, "")
$ws.Range("E3").Validation.InputMessage = "Entrez la tâche dans cette colonne sous ce titre."
$ws.Range("E3").Validation.ShowInput = $true
$ws.Range("E3").Validation.ShowError = $true

$ws.Range("F3").Validation.Add(0, 1, 1, "")
$ws.Range("F3").Validation.InputMessage = "Entrez les notes dans cette colonne sous ce titre."
$ws.Range("F3").Validation.ShowInput = $true
$ws.Range("F3").Validation.ShowError = $true

# ---- 9. Print setup (paper size / orientation) -------------------------
$ws.PageSetup.PaperSize = 9      # xlPaperA4
$ws.PageSetup.Orientation = 1    # xlPortrait

# ---- 10. Selection cursor, matching the saved file ---------------------
$ws.Range("B7").Select() | Out-Null

Write-Output "Rapport SEO sheet rebuilt."
